$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("F1").Value = "Vehicle Usage"
$ws.Range("G1").Value = "Assigned To"

# New data cells for rows 2-4
$ws.Range("F2").Value = "ASA"
$ws.Range("G2").Value = "Formal Name 1"

$ws.Range("F3").Value = "ASA"
$ws.Range("G3").Value = "Formal Name 2"

$ws.Range("F4").Value = "ASA"
$ws.Range("G4").Value = "Formal Name 3"

# Move the active selection, matching the author's final selection in the diff
$ws.Range("H6").Select()
